# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" timestamps (Overview: "Latest HO Xliff
# Generate Date") for the rows whose handoff xliff was just (re)generated, and
# marks those same rows' Priority as "ht" (handoff type) on the per-locale
# sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-29 02:22:56"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-29 02:22:51"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-29 02:22:56"
}
